$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$win.ScrollRow = 39
$ws.Range("B40").Select()
